$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '58.704.62'
$ws.Cells.Item(2, 5).Value = '  -3.94%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.574.05'
$ws.Cells.Item(3, 5).Value = '  -3.21%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '512.57'
$ws.Cells.Item(5, 5).Value = '  -3.74%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '147.24'
$ws.Cells.Item(6, 5).Value = '  -5.97%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.997'
$ws.Cells.Item(7, 5).Value = '  -0.08%  '

$ws.Cells.Item(8, 5).Value = '  -2.65%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '2.588.01'
$ws.Cells.Item(9, 5).Value = '  -3.28%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '6.33'
$ws.Cells.Item(10, 5).Value = '  -4.29%  '

$ws.Cells.Item(11, 5).Value = '  -4.96%  '

$ws.Cells.Item(12, 5).Value = '  -4.67%  '

$ws.Cells.Item(13, 5).Value = '  -0.76%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '3.022.18'
$ws.Cells.Item(14, 5).Value = '  -3.39%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '58.600.13'
$ws.Cells.Item(15, 5).Value = '  -4.08%  '

$ws.Cells.Item(16, 5).Value = '  -3.46%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.0000138'
$ws.Cells.Item(17, 5).Value = '  -4.16%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '2.582.88'
$ws.Cells.Item(18, 5).Value = '  -3.27%  '

$ws.Cells.Item(19, 2).Value = 'Polkadot'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '4.57'
$ws.Cells.Item(19, 5).Value = '  -4.38%  '

$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '348.11'
$ws.Cells.Item(20, 5).Value = '  -2.17%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '10.38'
$ws.Cells.Item(21, 5).Value = '  -3.18%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.08'
$ws.Cells.Item(22, 5).Value = '  -3.45%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  -0.01%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '60.86'
$ws.Cells.Item(24, 5).Value = '  -1.16%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.420'
$ws.Cells.Item(25, 5).Value = '  -3.03%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.998'
$ws.Cells.Item(26, 5).Value = '  -0.26%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.162'
$ws.Cells.Item(27, 5).Value = '  -4.20%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.685.76'
$ws.Cells.Item(28, 5).Value = '  -3.36%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.0₃0819'
$ws.Cells.Item(29, 5).Value = '  -5.26%  '

$ws.Cells.Item(30, 5).Value = '  -4.93%  '

$ws.Cells.Item(31, 5).Value = '  -0.07%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '6.07'
$ws.Cells.Item(32, 5).Value = '  -1.95%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '18.81'
$ws.Cells.Item(33, 5).Value = '  -3.99%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '149.95'
$ws.Cells.Item(34, 5).Value = '  +0.08%  '

$ws.Cells.Item(35, 5).Value = '  -4.78%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.943'
$ws.Cells.Item(36, 5).Value = '  +6.22%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.02'
$ws.Cells.Item(37, 5).Value = '  -3.29%  '

$ws.Cells.Item(38, 5).Value = '  -5.16%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.866'
$ws.Cells.Item(39, 5).Value = '  -5.90%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '36.17'
$ws.Cells.Item(40, 5).Value = '  -1.74%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '294.95'
$ws.Cells.Item(41, 5).Value = '  -3.92%  '

$ws.Cells.Item(42, 5).Value = '  -5.65%  '

$ws.Cells.Item(43, 5).Value = '  -6.00%  '

$ws.Cells.Item(44, 5).Value = '  -2.51%  '

$ws.Cells.Item(45, 5).Value = '  -0.13%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.611'
$ws.Cells.Item(46, 5).Value = '  -6.30%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0542'
$ws.Cells.Item(47, 5).Value = '  -4.38%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '19.25'
$ws.Cells.Item(48, 5).Value = '  -5.86%  '

$ws.Cells.Item(49, 5).Value = '  -3.85%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '4.72'
$ws.Cells.Item(50, 5).Value = '  -5.31%  '
